# Apply the target edit described by the diff:
#  - Append three new data rows (WishList data) to sheet1 ("testDataSheet"),
#    which also grows the shared strings table and the sheet dimension.
#  - Move the active cell selection on sheet1 from B8 to B14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: MyWishList1 / First WishList (note the trailing space, preserved as typed)
$ws.Cells.Item(4, 1).Value = "MyWishList1"
$ws.Cells.Item(4, 2).Value = "First WishList "

# Row 5: MyWishList2 / Second Wishlist
$ws.Cells.Item(5, 1).Value = "MyWishList2"
$ws.Cells.Item(5, 2).Value = "Second Wishlist"

# Row 6: MyWishList2 / Third Wishlist
$ws.Cells.Item(6, 1).Value = "MyWishList2"
$ws.Cells.Item(6, 2).Value = "Third Wishlist"

# Update the selected/active cell on the sheet to match the new target (B14)
$ws.Range("B14").Select() | Out-Null
